$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column A slightly so the exported OOXML width goes from 11 to 12.
# This engine's ColumnWidth -> OOXML width mapping adds ~0.8333 (5/6), so
# request 11.1666667 to land exactly on 12.
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666

# Update the data row values
$ws.Range("A2").Value = 0.93320580434120881
$ws.Range("B2").Value = 41511
$ws.Range("C2").Value = 976
$ws.Range("D2").Value = 976
